# Re-process the sheet with the newly curated dimensions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: metadata "type" identifiers
$ws.Range("A2").Value = "iaest-measure:horas-trabajadas"
$ws.Range("E2").Value = "sdmx-dimension:refArea"

# Row 3: dim/medida labels
$ws.Range("A3").Value = "medida"

# Row 4: datatype / URI kind
$ws.Range("A4").Value = "xsd:int"
$ws.Range("E4").Value = "URI-Comunidad"

# Row 5 (mapping file references) is no longer needed - remove it entirely.
$ws.Rows.Item(5).Delete()
